# Update countries & provincias Spain
# - Reorder "Polonia"/"Rumania" and "Emiratos Arabes Unidos" (country names)
#   so the table (sorted descending by "Casos totales") reflects the new
#   case counts.
# - Refresh the numeric stats for the affected rows plus Noruega (row 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 - Noruega: updated totals / new cases / active cases.
$ws.Cells.Item(20, 2).Value = 3970
$ws.Cells.Item(20, 3).Value = 199
$ws.Cells.Item(20, 5).Value = 3943

# Row 33 - now Polonia (was Rumania), with Polonia's refreshed stats.
$ws.Cells.Item(33, 1).Value = "Polonia"
$ws.Cells.Item(33, 2).Value = 1481
$ws.Cells.Item(33, 3).Value = 92
$ws.Cells.Item(33, 4).Value = 7
$ws.Cells.Item(33, 5).Value = 1457
$ws.Cells.Item(33, 6).Value = 3
$ws.Cells.Item(33, 7).Value = 1
$ws.Cells.Item(33, 8).Value = 17

# Row 34 - now Rumania (was Polonia), with Rumania's refreshed stats.
$ws.Cells.Item(34, 1).Value = "Rumania"
$ws.Cells.Item(34, 2).Value = 1452
$ws.Cells.Item(34, 3).Value = 160
$ws.Cells.Item(34, 4).Value = 139
$ws.Cells.Item(34, 5).Value = 1284
$ws.Cells.Item(34, 6).Value = 34
$ws.Cells.Item(34, 7).Value = 3
$ws.Cells.Item(34, 8).Value = 29

# Row 63 - now Emiratos Arabes Unidos (was Nueva Zelanda).
$ws.Cells.Item(63, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(63, 2).Value = 468
$ws.Cells.Item(63, 3).Value = 63
$ws.Cells.Item(63, 4).Value = 52
$ws.Cells.Item(63, 5).Value = 414
$ws.Cells.Item(63, 6).Value = 2
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 2

# Row 64 - now Nueva Zelanda (was Libano).
$ws.Cells.Item(64, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(64, 2).Value = 451
$ws.Cells.Item(64, 3).Value = 83
$ws.Cells.Item(64, 4).Value = 50
$ws.Cells.Item(64, 5).Value = 401
$ws.Cells.Item(64, 6).Value = 2
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0

# Row 65 - now Libano (was Argelia).
$ws.Cells.Item(65, 1).Value = "Libano"
$ws.Cells.Item(65, 2).Value = 412
$ws.Cells.Item(65, 3).Value = 21
$ws.Cells.Item(65, 4).Value = 27
$ws.Cells.Item(65, 5).Value = 377
$ws.Cells.Item(65, 6).Value = 3
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 8

# Row 66 - now Argelia (was Emiratos Arabes Unidos).
$ws.Cells.Item(66, 1).Value = "Argelia"
$ws.Cells.Item(66, 2).Value = 409
$ws.Cells.Item(66, 3).Value = 0
$ws.Cells.Item(66, 4).Value = 29
$ws.Cells.Item(66, 5).Value = 354
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = 26
